# 4.0.3 model and data
# Applies the CID "Check Input Data" workbook update:
#   - splits the single "trans/BVTQaZ/BVTQaZ.csv" row (Boolean sheet) into six
#     per-mode files (LDVs, HDVs, aircraft, rail, ships, motorbikes)
#   - splits the single "trans/VTQaZ/VTQaZ.csv" row (Boolean sheet) the same way
#   - adds a handful of trailing blank rows at the bottom of the Boolean sheet
#   - re-selects the "About" tab / updates the stored cell selections to match
#     where the author's cursor ended up after editing

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Boolean sheet: expand the two aggregate "...QaZ.csv" rows into six rows each
# ---------------------------------------------------------------------------
$boolSheet = $wb.Worksheets.Item("Boolean")

# --- trans/BVTQaZ/BVTQaZ.csv (currently row 17) -> six rows ---------------
[void]$boolSheet.Rows("17:21").Insert()
$bvtqaz = @(
    "trans/BVTQaZ/BVTQaZ-LDVs.csv",
    "trans/BVTQaZ/BVTQaZ-HDVs.csv",
    "trans/BVTQaZ/BVTQaZ-aircraft.csv",
    "trans/BVTQaZ/BVTQaZ-rail.csv",
    "trans/BVTQaZ/BVTQaZ-ships.csv",
    "trans/BVTQaZ/BVTQaZ-motorbikes.csv"
)
for ($i = 0; $i -lt $bvtqaz.Length; $i++) {
    $boolSheet.Cells.Item(17 + $i, 1).Value = $bvtqaz[$i]
}

# --- trans/VTQaZ/VTQaZ.csv (now shifted down to row 26) -> six rows -------
[void]$boolSheet.Rows("26:30").Insert()
$vtqaz = @(
    "trans/VTQaZ/VTQaZ-LDVs.csv",
    "trans/VTQaZ/VTQaZ-HDVs.csv",
    "trans/VTQaZ/VTQaZ-aircraft.csv",
    "trans/VTQaZ/VTQaZ-rail.csv",
    "trans/VTQaZ/VTQaZ-ships.csv",
    "trans/VTQaZ/VTQaZ-motorbikes.csv"
)
for ($i = 0; $i -lt $vtqaz.Length; $i++) {
    $boolSheet.Cells.Item(26 + $i, 1).Value = $vtqaz[$i]
}

# --- six trailing blank rows after the last data row (now row 32) ---------
[void]$boolSheet.Rows("33:38").Insert()

# Land the cursor/selection where the author left it in the commit
[void]$boolSheet.Range("A32").Select()

# ---------------------------------------------------------------------------
# Integer sheet: cursor moved to A13
# ---------------------------------------------------------------------------
$intSheet = $wb.Worksheets.Item("Integer")
[void]$intSheet.Range("A13").Select()

# ---------------------------------------------------------------------------
# About sheet becomes the active/selected tab again
# ---------------------------------------------------------------------------
$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Activate()
[void]$aboutSheet.Range("A1").Select()
